$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 21 (2025Q3) with refreshed metrics
$ws.Range("C21").Value = 157
$ws.Range("D21").Value = 144
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 41.26074498567336
